# Add 2022-Q4 data: insert a new quarter sheet + a new row in the summary sheet.
#
# NOTE: sheet/range references captured *before* a structural change (adding
# a sheet, inserting a row) can end up pointing at the wrong object once tab
# order / row numbers shift. So: do every structural change first, then
# re-fetch everything we need by name right before writing data into it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Structural changes
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# 1) Insert a new row 2 in "总计" for the 2022-Q4 summary line.
$summary.Rows.Item(2).Insert()

# 2) Insert a brand-new worksheet "2022-Q4" right after "总计" (it becomes the
#    2nd tab, pushing 2022-Q3 / 2022-Q2 / 2022-Q1 / 2021-Q4 one slot right).
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# "总计" - fill in the new row and renumber the index column.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# The inserted row inherits formatting from the row above (the bold/bordered
# header) - strip it back to plain so B2:D2 match the other data rows.
$summary.Range("B2:D2").ClearFormats()

# Give the new A2 the same "index column" look (bold / thin border /
# centered) as the other index cells below it.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 5
$summary.Range("D2").Value = 0.05

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# "2022-Q4" - header + 5 fund rows.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Item("2022-Q4")
$q3 = $wb.Worksheets.Item("2022-Q3")

# Reuse the header formatting (bold / border / centered) from the existing
# "2022-Q3" sheet's header row.
$q3.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Reuse the "index column" formatting for A2:A6.
$q3.Range("A2").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'010571"
$newSheet.Range("C2").Value = "新沃创新领航混合C"
$newSheet.Range("D2").Value = "'0.51"
$newSheet.Range("E2").Value = "'93.56"
$newSheet.Range("F2").Value = "'4.38"
$newSheet.Range("G2").Value = "'0.0223"
$newSheet.Range("H2").Value = 6

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'010570"
$newSheet.Range("C3").Value = "新沃创新领航混合A"
$newSheet.Range("D3").Value = "'0.24"
$newSheet.Range("E3").Value = "'93.56"
$newSheet.Range("F3").Value = "'4.38"
$newSheet.Range("G3").Value = "'0.0105"
$newSheet.Range("H3").Value = 6

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'012143"
$newSheet.Range("C4").Value = "新沃内需增长混合A"
$newSheet.Range("D4").Value = "'0.20"
$newSheet.Range("E4").Value = "'93.63"
$newSheet.Range("F4").Value = "'4.05"
$newSheet.Range("G4").Value = "'0.0081"
$newSheet.Range("H4").Value = 8

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'002564"
$newSheet.Range("C5").Value = "新沃通盈灵活配置混合"
$newSheet.Range("D5").Value = "'0.10"
$newSheet.Range("E5").Value = "'92.67"
$newSheet.Range("F5").Value = "'3.94"
$newSheet.Range("G5").Value = "'0.0039"
$newSheet.Range("H5").Value = 8

$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "'012144"
$newSheet.Range("C6").Value = "新沃内需增长混合C"
$newSheet.Range("D6").Value = "'0.04"
$newSheet.Range("E6").Value = "'93.63"
$newSheet.Range("F6").Value = "'4.05"
$newSheet.Range("G6").Value = "'0.0016"
$newSheet.Range("H6").Value = 8
